$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.508084416940065
$ws.Range("I2").Value = 0.552228264162119
$ws.Range("K2").Value = 0.382887811020303
$ws.Range("L2").Value = 0.502900789618378
$ws.Range("N2").Value = 0.484917310839545
